# OLX Monitor 2026-02-27 09:08 — append 9 new listing rows (243-251) to the
# "PODSUMOWANIE" sheet, continuing the existing per-listing detail table
# that runs from row 7 through row 242 (columns A-H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-formatting reference cells already present in the sheet, reused so the
# new rows land on the exact same style indexes as their "siblings":
#   style s=16 (bold, green "new today" highlight) -> F234 / F235
#   style s=14 (plain centered)                     -> F240
#   style s=15 (red "days listed" highlight)         -> F242
# Columns A/C/D/E keep the same style in every data row, so row 242 is a
# fine template for those regardless of which F-style a given row needs.

$newRows = @(
    @(243, "2026-02-27 09:08:20", "poqui", "Duży pokój z balkonem w 2pokojowym mieszkaniu blisko Politechniki", 1665, "25.02.2026", 1, "https://www.olx.pl/d/oferta/duzy-pokoj-z-balkonem-w-2pokojowym-mieszkaniu-blisko-politechniki-CID3-ID19xpQK.html", "duzy-pokoj-z-balkonem-w-2pokojowym-mieszkaniu-blisko-politechniki-CID3-ID19xpQK", 234),
    @(244, "2026-02-27 09:08:20", "poqui", "Nowoczesne mieszkanie 2-pokojowe z balkonem, blisko UMCS, KUL, UP", 2499, "25.02.2026", 1, "https://www.olx.pl/d/oferta/nowoczesne-mieszkanie-2-pokojowe-z-balkonem-blisko-umcs-kul-up-CID3-ID19xpwN.html", "nowoczesne-mieszkanie-2-pokojowe-z-balkonem-blisko-umcs-kul-up-CID3-ID19xpwN", 235),
    @(245, "2026-02-27 09:08:20", "poqui", "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy", 2499, "28.10.2025", 121, "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html", "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger", 242),
    @(246, "2026-02-27 09:08:20", "poqui", "Przytulny pokój blisko Politechniki – ul. Przytulna", 549, "10.10.2025", 139, "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html", "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz", 242),
    @(247, "2026-02-27 09:08:20", "poqui", "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza", 2049, "19.12.2025", 69, "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html", "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc", 242),
    @(248, "2026-02-27 09:08:20", "pokojewlublinie", "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12", 12640, "19.01.2026", 38, "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html", "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc", 240),
    @(249, "2026-02-27 09:08:20", "pokojewlublinie", "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58", 0, "11.08.2025", 199, "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html", "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm", 242),
    @(250, "2026-02-27 09:08:20", "dawnypatron", "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.", 730, "20.09.2024", 524, "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html", "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM", 242),
    @(251, "2026-02-27 09:08:20", "dawnypatron", "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14", 14690, "05.12.2025", 83, "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html", "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv", 242)
)

foreach ($r in $newRows) {
    $targetRow  = $r[0]
    $checkedAt  = $r[1]
    $profile    = $r[2]
    $title      = $r[3]
    $price      = $r[4]
    $datePosted = $r[5]
    $daysListed = $r[6]
    $url        = $r[7]
    $slug       = $r[8]
    $templateRow = $r[9]

    # Clone the whole template row's formatting (all 8 columns) onto the new
    # row first, so every cell lands on the correct existing style index.
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($templateRow, $c).Copy($ws.Cells.Item($targetRow, $c))
    }

    $ws.Cells.Item($targetRow, 1).Value = $checkedAt
    $ws.Cells.Item($targetRow, 2).Value = $profile
    $ws.Cells.Item($targetRow, 3).Value = $title
    $ws.Cells.Item($targetRow, 4).Value = $price

    # Column E holds a plain "DD.MM.YYYY" string in the source data. Excel's
    # usual auto-detection would silently reinterpret low-day values (day
    # <= 12, e.g. "10.10.2025") as a real date serial, which would change
    # both the stored type and the displayed text. Force the cell to Text
    # first so the literal string is preserved for every row, consistently.
    $eCell = $ws.Cells.Item($targetRow, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $datePosted

    $ws.Cells.Item($targetRow, 6).Value = $daysListed
    $ws.Cells.Item($targetRow, 7).Value = $url
    $ws.Cells.Item($targetRow, 8).Value = $slug
}
